# Daily attendance processing - 2026-01-05 06:50:44
# Reorders the "Recorded By" (column G) entries on the Session Analysis
# Results sheet: wherever the comma-separated list of recorders includes
# an exact "System" entry alongside other entries, the token order is
# reversed (System moves from front of the list to the back).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value()

    if ($null -eq $value) { continue }

    $parts = $value -split ","
    if ($parts.Count -gt 1) {
        $trimmed = @()
        foreach ($p in $parts) { $trimmed += $p.Trim() }

        $hasSystem = $false
        foreach ($t in $trimmed) {
            if ($t -ceq "System") { $hasSystem = $true }
        }

        if ($hasSystem) {
            $reversed = @()
            for ($i = $trimmed.Count - 1; $i -ge 0; $i--) {
                $reversed += $trimmed[$i]
            }
            $cell.Value = [string]::Join(", ", $reversed)
        }
    }
}
